$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain numeric cells (Qty executed upto date column)
$ws.Range("C8").Value = 50
$ws.Range("C9").Value = 95
$ws.Range("C10").Value = 65
$ws.Range("C11").Value = 26
$ws.Range("C12").Value = 56
$ws.Range("C13").Value = 28
$ws.Range("C14").Value = 58
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 70
$ws.Range("C17").Value = 26

# Amount cells stored as text (numeric-looking strings); use a leading
# apostrophe so Excel keeps them as text instead of coercing to numbers.
$ws.Range("G9").Value = "'24320.00"
$ws.Range("G10").Value = "'30680.00"
$ws.Range("G11").Value = "'17212.00"
$ws.Range("G13").Value = "'3808.00"
$ws.Range("G14").Value = "'1334.00"

$ws.Range("G19").Value = "'77354.00"
$ws.Range("H19").Value = "'77354.00"
$ws.Range("G21").Value = "'77354.00"
$ws.Range("H21").Value = "'77354.00"
